$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Global status text change: "Ready for handoff" -> "Handed back: in sync
#    with en-US". This string is shared by the Status column (C) on every
#    sheet, and by the zh-cn/de-de status columns (B/C) on the Overview sheet.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet ("Handed back" report for this locale):
#    - Latest Handback DateTime (H) gets a real timestamp
#    - Latest Target File (F) / Latest Handback File (G) are populated,
#      mirroring the Source File Name (A) / Latest Handoff File (D) entries
# ---------------------------------------------------------------------------
$wsZhCn.Range("H2").Value = "2016-03-19 20:49:24"
$wsZhCn.Range("H3").Value = "2016-03-19 20:49:24"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/240f8674e9139d12ebf961441c956dab8a77ce6f/e2e/7ba9af43-7a48-41ca-a957-fd186cc1e106.md",
    "",
    "",
    "7ba9af43-7a48-41ca-a957-fd186cc1e106.md")
$wsZhCn.Range("F2").Style = "HyperLink"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cb3305c97c83e4e899d43c74b93fd91bbcf60302/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7ba9af43-7a48-41ca-a957-fd186cc1e106.4b05dfabca316e48610534fc8cb5292a8223c94a.zh-cn.xlf",
    "",
    "",
    "7ba9af43-7a48-41ca-a957-fd186cc1e106.4b05dfabca316e48610534fc8cb5292a8223c94a.zh-cn.xlf")
$wsZhCn.Range("G2").Style = "HyperLink"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/240f8674e9139d12ebf961441c956dab8a77ce6f/e2e/b25aff4e-4002-4458-920d-977b48fe4580.md",
    "",
    "",
    "b25aff4e-4002-4458-920d-977b48fe4580.md")
$wsZhCn.Range("F3").Style = "HyperLink"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cb3305c97c83e4e899d43c74b93fd91bbcf60302/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b25aff4e-4002-4458-920d-977b48fe4580.c62c568c1820b6abc3d85b4d5201f524377df177.zh-cn.xlf",
    "",
    "",
    "b25aff4e-4002-4458-920d-977b48fe4580.c62c568c1820b6abc3d85b4d5201f524377df177.zh-cn.xlf")
$wsZhCn.Range("G3").Style = "HyperLink"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same shape of update, different timestamp/locale targets.
# ---------------------------------------------------------------------------
$wsDeDe.Range("H2").Value = "2016-03-19 20:49:30"
$wsDeDe.Range("H3").Value = "2016-03-19 20:49:30"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/240f8674e9139d12ebf961441c956dab8a77ce6f/e2e/7ba9af43-7a48-41ca-a957-fd186cc1e106.md",
    "",
    "",
    "7ba9af43-7a48-41ca-a957-fd186cc1e106.md")
$wsDeDe.Range("F2").Style = "HyperLink"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e05045180c7ef3f0894367afdb7fdf6722ff2d9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7ba9af43-7a48-41ca-a957-fd186cc1e106.4b05dfabca316e48610534fc8cb5292a8223c94a.de-de.xlf",
    "",
    "",
    "7ba9af43-7a48-41ca-a957-fd186cc1e106.4b05dfabca316e48610534fc8cb5292a8223c94a.de-de.xlf")
$wsDeDe.Range("G2").Style = "HyperLink"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/240f8674e9139d12ebf961441c956dab8a77ce6f/e2e/b25aff4e-4002-4458-920d-977b48fe4580.md",
    "",
    "",
    "b25aff4e-4002-4458-920d-977b48fe4580.md")
$wsDeDe.Range("F3").Style = "HyperLink"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e05045180c7ef3f0894367afdb7fdf6722ff2d9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b25aff4e-4002-4458-920d-977b48fe4580.c62c568c1820b6abc3d85b4d5201f524377df177.de-de.xlf",
    "",
    "",
    "b25aff4e-4002-4458-920d-977b48fe4580.c62c568c1820b6abc3d85b4d5201f524377df177.de-de.xlf")
$wsDeDe.Range("G3").Style = "HyperLink"
